$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Angkatan" column (column C),
# shifting it to column D, then fill in the new "Prodi ID" column.
$ws.Columns.Item(3).Insert()

$ws.Cells.Item(1, 3).Value = "Prodi ID"
$ws.Cells.Item(2, 3).Value = 13201
$ws.Cells.Item(3, 3).Value = 14001

# Match the styling applied to the new data cell (Arial, dark grey).
$c = $ws.Cells.Item(2, 3)
$c.Font.Name = "Arial"
$c.Font.Color = 3355443

# Narrow the new column a touch so it doesn't look oversized.
$ws.Columns.Item(3).ColumnWidth = 7.29

# Page setup: portrait orientation, 300dpi print quality.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = 300

# Restore the final selected cell as recorded by Excel on save.
$null = $ws.Range("D11").Select()
